# #5: property boat&car done
#
# The "汽車" (car) sheet only ever had a single data row whose "header" row
# (row 1) was a stray duplicate of the data instead of real column titles,
# and it was missing the trailing metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) that
# every other property sheet ("土地"/land, "建物"/building, "股票"/stock, ...)
# already carries. This brings "汽車" up to the same schema, introducing a
# new "capacity" field (engine displacement) in place of "area".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")
$landSheet = $wb.Worksheets.Item("土地")

# Pull over the bold/bordered header formatting for the newly added
# H:N columns from the already-complete "land" sheet (same visual style
# used for B1:G1 already).
$landSheet.Range("H1:N1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Row 1: real column headers (previously this row just repeated the
#     data row's values) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: extend the existing car record (name/year/owner/register
#     date/reason/value stay as-is) with the shared metadata columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "date" must stay literal text ("2012-04-19"), not get auto-coerced into
# a serial date value, so copy it straight from the land sheet's already-
# correct, identically-formatted "date" cell instead of retyping it.
$landSheet.Range("K2").Copy()
$ws.Range("J2").PasteSpecial(-4104)

$ws.Range("K2").Value = "賴士葆"
$ws.Range("L2").Value = 866
$ws.Range("M2").Value = "tmp9edb1"
$ws.Range("N2").Value = 32
